$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 15
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H15").Value = 137268.06
$ws_ALC.Range("I15").Value = 137268.06
$ws_ALC.Range("K15").Value = 411804.18
$ws_ALC.Range("M15").Value = -411635.18

# Sheet ALC, Row 112
$ws_ALC.Range("H112").Value = 5803796.5
$ws_ALC.Range("J112").Value = 5803796.5
$ws_ALC.Range("L112").Value = 17411389.5
$ws_ALC.Range("N112").Value = -17413605.5

# Sheet ALC, Row 138
$ws_ALC.Range("H138").Value = 9526023
$ws_ALC.Range("I138").Value = 1852671
$ws_ALC.Range("J138").Value = 13160769
$ws_ALC.Range("K138").Value = 5558013
$ws_ALC.Range("L138").Value = 39482307
$ws_ALC.Range("M138").Value = -5552873
$ws_ALC.Range("N138").Value = -39492587

# Sheet ALC, Row 141
$ws_ALC.Range("H141").Value = 2529.5454
$ws_ALC.Range("I141").Value = 2203.889
$ws_ALC.Range("J141").Value = 3995
$ws_ALC.Range("K141").Value = 6611.667
$ws_ALC.Range("L141").Value = 11985
$ws_ALC.Range("M141").Value = -1431.667
$ws_ALC.Range("N141").Value = -22345

# Sheet ARM, Row 2
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 1118.5333
$ws_ARM.Range("I2").Value = 700.8
$ws_ARM.Range("J2").Value = 1954
$ws_ARM.Range("K2").Value = 700.8
$ws_ARM.Range("L2").Value = 1954
$ws_ARM.Range("M2").Value = -587.8
$ws_ARM.Range("N2").Value = -2180

# Sheet ARM, Row 74
$ws_ARM.Range("H74").Value = 8209.789000000001
$ws_ARM.Range("I74").Value = 2203
$ws_ARM.Range("J74").Value = 18507.143
$ws_ARM.Range("K74").Value = 2203
$ws_ARM.Range("L74").Value = 18507.143
$ws_ARM.Range("M74").Value = -1329
$ws_ARM.Range("N74").Value = -20255.143

# Sheet ARM, Row 77
$ws_ARM.Range("H77").Value = 8209.789000000001
$ws_ARM.Range("I77").Value = 2203
$ws_ARM.Range("J77").Value = 18507.143
$ws_ARM.Range("K77").Value = 11015
$ws_ARM.Range("L77").Value = 92535.715
$ws_ARM.Range("M77").Value = -6647
$ws_ARM.Range("N77").Value = -101271.715

# Sheet ARM, Row 110
$ws_ARM.Range("H110").Value = 1104.5454
$ws_ARM.Range("I110").Value = 905.55554
$ws_ARM.Range("J110").Value = 2000
$ws_ARM.Range("K110").Value = 905.55554
$ws_ARM.Range("L110").Value = 2000
$ws_ARM.Range("M110").Value = 1139.44446
$ws_ARM.Range("N110").Value = -6090

# Sheet ARM, Row 116
$ws_ARM.Range("H116").Value = 1118.5333
$ws_ARM.Range("I116").Value = 700.8
$ws_ARM.Range("J116").Value = 1954
$ws_ARM.Range("K116").Value = 700.8
$ws_ARM.Range("L116").Value = 1954
$ws_ARM.Range("M116").Value = 1593.2
$ws_ARM.Range("N116").Value = -6542

# Sheet ARM, Row 132
$ws_ARM.Range("H132").Value = 3941.8462
$ws_ARM.Range("I132").Value = 3760.889
$ws_ARM.Range("K132").Value = 11282.667
$ws_ARM.Range("M132").Value = -8752.667000000001

# Sheet BSM, Row 3
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 1118.5333
$ws_BSM.Range("I3").Value = 700.8
$ws_BSM.Range("J3").Value = 1954
$ws_BSM.Range("K3").Value = 700.8
$ws_BSM.Range("L3").Value = 1954
$ws_BSM.Range("M3").Value = -586.8
$ws_BSM.Range("N3").Value = -2182

# Sheet BSM, Row 94
$ws_BSM.Range("H94").Value = 1895.125
$ws_BSM.Range("I94").Value = 1994.4286
$ws_BSM.Range("J94").Value = 1200
$ws_BSM.Range("K94").Value = 1994.4286
$ws_BSM.Range("L94").Value = 1200
$ws_BSM.Range("M94").Value = -1543.4286
$ws_BSM.Range("N94").Value = -2102

# Sheet BSM, Row 99
$ws_BSM.Range("H99").Value = 793.3570999999999
$ws_BSM.Range("I99").Value = 793.3570999999999
$ws_BSM.Range("K99").Value = 793.3570999999999
$ws_BSM.Range("M99").Value = 704.6429000000001

# Sheet BSM, Row 107
$ws_BSM.Range("H107").Value = 431.31818
$ws_BSM.Range("I107").Value = 276.77777
$ws_BSM.Range("K107").Value = 276.77777
$ws_BSM.Range("M107").Value = 1643.22223

# Sheet BSM, Row 134
$ws_BSM.Range("H134").Value = 5236.5713
$ws_BSM.Range("I134").Value = 4516.75
$ws_BSM.Range("K134").Value = 13550.25
$ws_BSM.Range("M134").Value = -11015.25

# Sheet CRP, Row 16
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 672.8570999999999
$ws_CRP.Range("I16").Value = 442.2
$ws_CRP.Range("K16").Value = 442.2
$ws_CRP.Range("M16").Value = -155.2

# Sheet CRP, Row 58
$ws_CRP.Range("H58").Value = 2110.875
$ws_CRP.Range("I58").Value = 1392
$ws_CRP.Range("K58").Value = 1392
$ws_CRP.Range("M58").Value = -1189

# Sheet CRP, Row 113
$ws_CRP.Range("H113").Value = 672.8570999999999
$ws_CRP.Range("I113").Value = 442.2
$ws_CRP.Range("K113").Value = 442.2
$ws_CRP.Range("M113").Value = 1727.8

# Sheet CRP, Row 136
$ws_CRP.Range("H136").Value = 2110.875
$ws_CRP.Range("I136").Value = 1392
$ws_CRP.Range("K136").Value = 4176
$ws_CRP.Range("M136").Value = -1626

# Sheet CUL, Row 12
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H12").Value = 31.588236
$ws_CUL.Range("J12").Value = 128.33333
$ws_CUL.Range("L12").Value = 384.99999
$ws_CUL.Range("N12").Value = -730.99999

# Sheet GSM, Row 107
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H107").Value = 167.04546
$ws_GSM.Range("I107").Value = 170.36363
$ws_GSM.Range("J107").Value = 163.72728
$ws_GSM.Range("K107").Value = 170.36363
$ws_GSM.Range("L107").Value = 163.72728
$ws_GSM.Range("M107").Value = 1749.63637
$ws_GSM.Range("N107").Value = -4003.72728

# Sheet GSM, Row 113
$ws_GSM.Range("H113").Value = 0
$ws_GSM.Range("I113").Value = 0
$ws_GSM.Range("J113").Value = 0
$ws_GSM.Range("K113").Value = 0
$ws_GSM.Range("L113").Value = 0
$ws_GSM.Range("M113").ClearContents()
$ws_GSM.Range("N113").ClearContents()

# Sheet LTW, Row 7
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2849.32
$ws_LTW.Range("I7").Value = 1747.5714
$ws_LTW.Range("K7").Value = 1747.5714
$ws_LTW.Range("M7").Value = -1635.5714

# Sheet LTW, Row 46
$ws_LTW.Range("H46").Value = 456.94116
$ws_LTW.Range("I46").Value = 436.66666
$ws_LTW.Range("J46").Value = 479.75
$ws_LTW.Range("K46").Value = 436.66666
$ws_LTW.Range("L46").Value = 479.75
$ws_LTW.Range("M46").Value = -248.66666
$ws_LTW.Range("N46").Value = -855.75

# Sheet LTW, Row 61
$ws_LTW.Range("H61").Value = 8313
$ws_LTW.Range("I61").Value = 8546.182000000001
$ws_LTW.Range("J61").Value = 7800
$ws_LTW.Range("K61").Value = 8546.182000000001
$ws_LTW.Range("L61").Value = 7800
$ws_LTW.Range("M61").Value = -8344.182000000001
$ws_LTW.Range("N61").Value = -8204

# Sheet LTW, Row 93
$ws_LTW.Range("H93").Value = 642.6799999999999
$ws_LTW.Range("I93").Value = 628.9
$ws_LTW.Range("J93").Value = 697.8
$ws_LTW.Range("K93").Value = 628.9
$ws_LTW.Range("L93").Value = 697.8
$ws_LTW.Range("M93").Value = 619.1
$ws_LTW.Range("N93").Value = -3193.8

# Sheet LTW, Row 108
$ws_LTW.Range("H108").Value = 30000
$ws_LTW.Range("J108").Value = 30000
$ws_LTW.Range("L108").Value = 30000
$ws_LTW.Range("N108").Value = -37680

# Sheet LTW, Row 113
$ws_LTW.Range("H113").Value = 8313
$ws_LTW.Range("I113").Value = 8546.182000000001
$ws_LTW.Range("J113").Value = 7800
$ws_LTW.Range("K113").Value = 8546.182000000001
$ws_LTW.Range("L113").Value = 7800
$ws_LTW.Range("M113").Value = -6376.182000000001
$ws_LTW.Range("N113").Value = -12140

# Sheet LTW, Row 126
$ws_LTW.Range("H126").Value = 2849.32
$ws_LTW.Range("I126").Value = 1747.5714
$ws_LTW.Range("K126").Value = 5242.7142
$ws_LTW.Range("M126").Value = -2772.7142

# Sheet LTW, Row 136
$ws_LTW.Range("H136").Value = 4657.4707
$ws_LTW.Range("I136").Value = 2995.1
$ws_LTW.Range("J136").Value = 7032.2856
$ws_LTW.Range("K136").Value = 8985.299999999999
$ws_LTW.Range("L136").Value = 21096.8568
$ws_LTW.Range("M136").Value = -6435.299999999999
$ws_LTW.Range("N136").Value = -26196.8568

# Sheet WVR, Row 113
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H113").Value = 229.66667
$ws_WVR.Range("I113").Value = 195
$ws_WVR.Range("J113").Value = 299
$ws_WVR.Range("K113").Value = 585
$ws_WVR.Range("L113").Value = 897
$ws_WVR.Range("M113").Value = 1585
$ws_WVR.Range("N113").Value = -5237

# Sheet WVR, Row 126
$ws_WVR.Range("H126").Value = 41383.52
$ws_WVR.Range("I126").Value = 44843.176
$ws_WVR.Range("J126").Value = 1597.5
$ws_WVR.Range("K126").Value = 134529.528
$ws_WVR.Range("L126").Value = 4792.5
$ws_WVR.Range("M126").Value = -132059.528
$ws_WVR.Range("N126").Value = -9732.5
